# add userform screenshots, reorder consent to match appt tracking sheet
#
# The "consent" column (previously the last data column, Y) is moved to
# become the first of the trailing block of columns (U), pushing
# device/video/sound/fun one column to the right (U->V, V->W, W->X, X->Y).
# Three new participant rows (6-8) are appended, and C2's test date is
# corrected from 12/15 to 12/14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Header row: re-point U1:Y1 so "consent" leads, then device/video/
#    sound/fun follow (was device/video/sound/fun/consent).
# ---------------------------------------------------------------------
$ws.Range("U1").Value = "consent"
$ws.Range("V1").Value = "device"
$ws.Range("W1").Value = "video"
$ws.Range("X1").Value = "sound"
$ws.Range("Y1").Value = "fun"

# ---------------------------------------------------------------------
# 2. Fix the test_date typo on row 2 (was 12/15/2020, should be 12/14).
# ---------------------------------------------------------------------
$ws.Range("C2").Value = 44179

# ---------------------------------------------------------------------
# 3. Shift the existing per-row consent/device/video/sound/fun values
#    the same way as the header for the four existing data rows
#    (2-5): old U (device) -> V, old V -> W, old W -> X, old X -> Y,
#    old Y (consent) -> U. Capture the "old" values before overwriting.
# ---------------------------------------------------------------------
foreach ($r in 2..5) {
    $oldU = $ws.Range("U$r").Value2
    $oldV = $ws.Range("V$r").Value2
    $oldW = $ws.Range("W$r").Value2
    $oldX = $ws.Range("X$r").Value2
    $oldY = $ws.Range("Y$r").Value2

    $ws.Range("U$r").Value = $oldY
    $ws.Range("V$r").Value = $oldU
    $ws.Range("W$r").Value = $oldV
    $ws.Range("X$r").Value = $oldW
    $ws.Range("Y$r").Value = $oldX
}

# ---------------------------------------------------------------------
# 4. Append three new participant rows (userform screenshots added for
#    these new sessions).
# ---------------------------------------------------------------------

# Row 6 - participant 5, studyname_20201216_01
$ws.Cells.Item(2, 3).Copy($ws.Cells.Item(6, 3)) | Out-Null
$ws.Range("A6").Value = "studyname_20201216_01"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 44181
$ws.Range("D6").Value = "MZ"
$ws.Range("E6").Value = "Zoom"
$ws.Range("F6").Value = "fill manually from REDCap"
$ws.Range("G6").Value = "fill manually from REDCap"
$ws.Range("H6").Value = "fill manually from REDCap"
$ws.Range("I6").Value = "fill manually from REDCap"
$ws.Range("J6").Value = "condition3"
$ws.Range("K6").Value = "left"
$ws.Range("L6").Value = "Yes"
$ws.Range("M6").Value = "blueberries"
$ws.Range("N6").Value = "blue"
$ws.Range("O6").Value = "left"
$ws.Range("P6").Value = "right"
$ws.Range("Q6").Value = "bottom"
$ws.Range("R6").Value = "hard"
$ws.Range("S6").Value = "a little hard"
$ws.Range("T6").Value = "eeeeeeeeee"
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = "laptop"
$ws.Range("W6").Value = 2
$ws.Range("X6").Value = 3
$ws.Range("Y6").Value = 5

# Row 7 - participant 6, studyname_20201216_02
$ws.Cells.Item(2, 3).Copy($ws.Cells.Item(7, 3)) | Out-Null
$ws.Range("A7").Value = "studyname_20201216_02"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 44181
$ws.Range("D7").Value = "MZ"
$ws.Range("E7").Value = "Zoom"
$ws.Range("F7").Value = "fill manually from REDCap"
$ws.Range("G7").Value = "fill manually from REDCap"
$ws.Range("H7").Value = "fill manually from REDCap"
$ws.Range("I7").Value = "fill manually from REDCap"
$ws.Range("J7").Value = "condition1"
$ws.Range("K7").Value = "left"
$ws.Range("L7").Value = "Yes"
$ws.Range("M7").Value = "blueberries"
$ws.Range("N7").Value = "blue"
$ws.Range("O7").Value = "left"
$ws.Range("P7").Value = "right"
$ws.Range("Q7").Value = "bottom"
$ws.Range("R7").Value = "easy"
$ws.Range("S7").Value = "a little easy"
$ws.Range("T7").Value = "jumped around slides in a random order!"
$ws.Range("U7").Value = 1
$ws.Range("V7").Value = "tablet"
$ws.Range("W7").Value = 1
$ws.Range("X7").Value = 4
$ws.Range("Y7").Value = 5

# Row 8 - participant 7, studyname_20201217_01
$ws.Cells.Item(2, 3).Copy($ws.Cells.Item(8, 3)) | Out-Null
$ws.Range("A8").Value = "studyname_20201217_01"
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 44182
$ws.Range("D8").Value = "MZ"
$ws.Range("E8").Value = "Zoom"
$ws.Range("F8").Value = "fill manually from REDCap"
$ws.Range("G8").Value = "fill manually from REDCap"
$ws.Range("H8").Value = "fill manually from REDCap"
$ws.Range("I8").Value = "fill manually from REDCap"
$ws.Range("J8").Value = "condition1"
$ws.Range("K8").Value = "left"
$ws.Range("L8").Value = "Yes"
$ws.Range("M8").Value = "blueberries"
$ws.Range("N8").Value = "blue"
$ws.Range("O8").Value = "left"
$ws.Range("P8").Value = "left"
$ws.Range("Q8").Value = "bottom"
$ws.Range("R8").Value = "easy"
$ws.Range("S8").Value = "a little easy"
$ws.Range("T8").Value = "blah"
$ws.Range("U8").Value = 3
$ws.Range("V8").Value = "laptop"
$ws.Range("W8").Value = 2
$ws.Range("X8").Value = 3
$ws.Range("Y8").Value = 4

# ---------------------------------------------------------------------
# 5. Column width tweaks (new/changed columns to accommodate the new
#    data) and refreshed selection/scroll position.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 9.5
$ws.Columns.Item(3).ColumnWidth = 10.333333333333334
$ws.Columns.Item(6).ColumnWidth = 9.833333333333334
$ws.Columns.Item(9).ColumnWidth = 13.666666666666666
$ws.Columns.Item(13).ColumnWidth = 10.833333333333334
$ws.Columns.Item(18).ColumnWidth = 10.5

$ws.Range("V9").Select()
